# Adds the new "Entry" WY group (rows 15-20) to Sheet1 and a summary
# AVERAGE formula in P22, matching the upstream commit that appended the
# 2025 "Entry" survival-rate rows to HatcheryWinterRunSurvival.xlsx.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Data rows 15-20 ("Entry" metric), columns A-H
# ---------------------------------------------------------------------

# -- Row 15 --------------------------------------------------------
$ws.Range("A15").Value = "Entry"
$ws.Range("B15").Formula = "=C15-1"
$ws.Range("C15").Value = 2020
$ws.Range("D15").Value = 13.2
$ws.Range("E15").Value = 1.5
$ws.Range("F15").Value = 10.5
$ws.Range("G15").Value = 16.5
$ws.Range("H15").Value = 93.4

# -- Row 16 ----------------------------------------------------------
$ws.Range("A16").Value = "Entry"
$ws.Range("B16:B19").Formula = "=C16-1"
$ws.Range("C16").Value = 2021
$ws.Range("D16").Value = 10.1
$ws.Range("E16").Value = 1.3
$ws.Range("F16").Value = 7.8
$ws.Range("G16").Value = 12.9
$ws.Range("H16").Value = 100

# -- Row 17 ------------------------------------------------------------
$ws.Range("A17").Value = "Entry"
$ws.Range("C17").Value = 2022
$ws.Range("D17").Value = 13.4
$ws.Range("E17").Value = 1.4
$ws.Range("F17").Value = 10.8
$ws.Range("G17").Value = 16.4
$ws.Range("H17").Value = 100

# -- Row 18 ------------------------------------------------------------
$ws.Range("A18").Value = "Entry"
$ws.Range("C18").Value = 2023
$ws.Range("D18").Value = 14
$ws.Range("E18").Value = 1.3
$ws.Range("F18").Value = 11.7
$ws.Range("G18").Value = 16.8
$ws.Range("H18").Value = 64.8

# -- Row 19 ------------------------------------------------------------
$ws.Range("A19").Value = "Entry"
$ws.Range("C19").Value = 2024
$ws.Range("D19").Value = 35.3
$ws.Range("E19").Value = 2.1
$ws.Range("F19").Value = 31.3
$ws.Range("G19").Value = 39.5
$ws.Range("H19").Value = 76

# -- Row 20 (plain values, no formula for B20) --------------------------
$ws.Range("A20").Value = "Entry"
$ws.Range("B20").Value = 2024
$ws.Range("C20").Value = 2025
$ws.Range("D20").Value = 38
$ws.Range("E20").Value = 1.8
$ws.Range("F20").Value = 34.4
$ws.Range("G20").Value = 41.6
$ws.Range("H20").Value = 70

# ---------------------------------------------------------------------
# 2. Formatting: mirror the look of the "Benicia" block (rows 8-14) --
#    column A = bold-less label font, B/C = year font with full box
#    border, D = dark-grey font, E-H = normal font - all boxed.
# ---------------------------------------------------------------------

$ws.Range("A15:A20").Font.Name = "Arial"
$ws.Range("A15:A20").Font.Size = 12

$ws.Range("B15:C20").Font.Name = "Arial"
$ws.Range("B15:C20").Font.Size = 12
$ws.Range("B15:C20").Borders.LineStyle = 1
$ws.Range("B15:C20").Borders.Weight = 2

$ws.Range("D15:D19").Font.Name = "Arial"
$ws.Range("D15:D19").Font.Size = 12
$ws.Range("D15:D19").Font.Color = 3355443

$ws.Range("E15:H19").Font.Name = "Arial"
$ws.Range("E15:H19").Font.Size = 12

# Vertical-only (left+right) thin borders for D15:H19, matching the
# "Entry" block's box-less look (no top/bottom rules between rows).
$rngLR = $ws.Range("D15:H15,F16:H19,D17:D19,E17:E19")
$rngLR.Borders.Item(7).LineStyle = 1
$rngLR.Borders.Item(7).Weight = 2
$rngLR.Borders.Item(10).LineStyle = 1
$rngLR.Borders.Item(10).Weight = 2

# Row 15 D:H and row 16 F:H keep both left+right edges; rows 17-19 only
# show the inner D|E divider (D = left edge, E = right edge).
$ws.Range("D17:D19").Borders.Item(10).LineStyle = -4142
$ws.Range("E17:E19").Borders.Item(7).LineStyle = -4142

# Row 20 and D16:E16 stay on the default (no border, default font) style,
# so nothing further to do there.

# ---------------------------------------------------------------------
# 3. Row heights to match the rest of the table
# ---------------------------------------------------------------------
$ws.Range("A15:P20").RowHeight = 15.75

# ---------------------------------------------------------------------
# 4. Summary formula
# ---------------------------------------------------------------------
$ws.Range("P22").Formula = "=AVERAGE(D15:D20)"

# ---------------------------------------------------------------------
# 5. Selection / view state to match saved file
# ---------------------------------------------------------------------
$ws.Range("P23").Select()
